$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1146.8182
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1146.8182
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 3440.4546
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -6934.4546
$ws.Range("H129").Value = 912.8913
$ws.Range("I129").Value = 379.36365
$ws.Range("J129").Value = 1080.5714
$ws.Range("K129").Value = 1138.09095
$ws.Range("L129").Value = 3241.7142
$ws.Range("M129").Value = 3861.90905
$ws.Range("N129").Value = -13241.7142
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 206127.67
$ws.Range("I74").Value = 1998.2
$ws.Range("J74").Value = 1113369.8
$ws.Range("K74").Value = 1998.2
$ws.Range("L74").Value = 1113369.8
$ws.Range("M74").Value = -1124.2
$ws.Range("N74").Value = -1115117.8
$ws.Range("H77").Value = 206127.67
$ws.Range("I77").Value = 1998.2
$ws.Range("J77").Value = 1113369.8
$ws.Range("K77").Value = 9991
$ws.Range("L77").Value = 5566849
$ws.Range("M77").Value = -5623
$ws.Range("N77").Value = -5575585
$ws.Range("H110").Value = 2266.95
$ws.Range("I110").Value = 2438.5
$ws.Range("J110").Value = 1866.6666
$ws.Range("K110").Value = 2438.5
$ws.Range("L110").Value = 1866.6666
$ws.Range("M110").Value = -393.5
$ws.Range("N110").Value = -5956.6666
$ws.Range("H112").Value = 33129
$ws.Range("I112").Value = 40000
$ws.Range("K112").Value = 40000
$ws.Range("M112").Value = -38523
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 28780
$ws.Range("I87").Value = 20000
$ws.Range("J87").Value = 30975
$ws.Range("K87").Value = 20000
$ws.Range("L87").Value = 30975
$ws.Range("M87").Value = -18752
$ws.Range("N87").Value = -33471
$ws.Range("H90").Value = 28780
$ws.Range("I90").Value = 20000
$ws.Range("J90").Value = 30975
$ws.Range("K90").Value = 60000
$ws.Range("L90").Value = 92925
$ws.Range("M90").Value = -53760
$ws.Range("N90").Value = -105405
$ws.Range("H99").Value = 1747.4231
$ws.Range("I99").Value = 1499.9286
$ws.Range("J99").Value = 2036.1666
$ws.Range("K99").Value = 1499.9286
$ws.Range("L99").Value = 2036.1666
$ws.Range("M99").Value = -1.92859999999996
$ws.Range("N99").Value = -5032.1666
$ws.Range("H105").Value = 1704.1
$ws.Range("I105").Value = 1259.3334
$ws.Range("K105").Value = 1259.3334
$ws.Range("M105").Value = 487.6666
$ws.Range("H110").Value = 34180
$ws.Range("J110").Value = 34180
$ws.Range("L110").Value = 34180
$ws.Range("N110").Value = -42360
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2012.836
$ws.Range("I31").Value = 1220.1154
$ws.Range("J31").Value = 2601.7144
$ws.Range("K31").Value = 1220.1154
$ws.Range("L31").Value = 2601.7144
$ws.Range("M31").Value = -925.1153999999999
$ws.Range("N31").Value = -3191.7144
$ws.Range("H34").Value = 2012.836
$ws.Range("I34").Value = 1220.1154
$ws.Range("J34").Value = 2601.7144
$ws.Range("K34").Value = 1220.1154
$ws.Range("L34").Value = 2601.7144
$ws.Range("M34").Value = -1018.1154
$ws.Range("N34").Value = -3005.7144
$ws.Range("H86").Value = 296317.88
$ws.Range("I86").Value = 386638.78
$ws.Range("J86").Value = 2775
$ws.Range("K86").Value = 386638.78
$ws.Range("L86").Value = 2775
$ws.Range("M86").Value = -385515.78
$ws.Range("N86").Value = -5021
$ws.Range("H89").Value = 296317.88
$ws.Range("I89").Value = 386638.78
$ws.Range("J89").Value = 2775
$ws.Range("K89").Value = 1933193.9
$ws.Range("L89").Value = 13875
$ws.Range("M89").Value = -1927577.9
$ws.Range("N89").Value = -25107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 720.2692
$ws.Range("J5").Value = 1015.8461
$ws.Range("L5").Value = 3047.5383
$ws.Range("N5").Value = -3271.5383
$ws.Range("H68").Value = 1423.125
$ws.Range("I68").Value = 998
$ws.Range("K68").Value = 2994
$ws.Range("M68").Value = -2183
$ws.Range("H71").Value = 1423.125
$ws.Range("I71").Value = 998
$ws.Range("K71").Value = 8982
$ws.Range("M71").Value = -4926
$ws.Range("H135").Value = 720.2692
$ws.Range("J135").Value = 1015.8461
$ws.Range("L135").Value = 9142.6149
$ws.Range("N135").Value = -14212.6149
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 59248136
$ws.Range("I70").Value = 276470900
$ws.Range("J70").Value = 5560.5454
$ws.Range("K70").Value = 276470900
$ws.Range("L70").Value = 5560.5454
$ws.Range("M70").Value = -276470630
$ws.Range("N70").Value = -6100.5454
$ws.Range("H73").Value = 59248136
$ws.Range("I73").Value = 276470900
$ws.Range("J73").Value = 5560.5454
$ws.Range("K73").Value = 276470900
$ws.Range("L73").Value = 5560.5454
$ws.Range("M73").Value = -276469964
$ws.Range("N73").Value = -7432.5454
$ws.Range("H111").Value = 22764.334
$ws.Range("J111").Value = 22764.334
$ws.Range("L111").Value = 22764.334
$ws.Range("N111").Value = -28898.334
$ws.Range("H126").Value = 3047.3809
$ws.Range("I126").Value = 1899.6428
$ws.Range("J126").Value = 5342.857
$ws.Range("K126").Value = 5698.928400000001
$ws.Range("L126").Value = 16028.571
$ws.Range("M126").Value = -3228.928400000001
$ws.Range("N126").Value = -20968.571
$ws.Range("H132").Value = 2464.7878
$ws.Range("I132").Value = 2001.8
$ws.Range("K132").Value = 6005.4
$ws.Range("M132").Value = -3475.4
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 422.0625
$ws.Range("I16").Value = 268.45456
$ws.Range("K16").Value = 268.45456
$ws.Range("M16").Value = -98.45456000000001
$ws.Range("H92").Value = 27950
$ws.Range("J92").Value = 27950
$ws.Range("L92").Value = 27950
$ws.Range("N92").Value = -32942
$ws.Range("H110").Value = 19000
$ws.Range("J110").Value = 19000
$ws.Range("L110").Value = 19000
$ws.Range("N110").Value = -27180
$ws.Range("H132").Value = 11911813
$ws.Range("I132").Value = 16674920
$ws.Range("K132").Value = 50024760
$ws.Range("M132").Value = -50022230
$ws.Range("H136").Value = 5001.282
$ws.Range("I136").Value = 6967.7144
$ws.Range("J136").Value = 2707.111
$ws.Range("K136").Value = 20903.1432
$ws.Range("L136").Value = 8121.333
$ws.Range("M136").Value = -18353.1432
$ws.Range("N136").Value = -13221.333
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 22006.178
$ws.Range("I100").Value = 60027.3
$ws.Range("K100").Value = 120054.6
$ws.Range("M100").Value = -119513.6
$ws.Range("H132").Value = 1934.25
$ws.Range("I132").Value = 1000.8
$ws.Range("J132").Value = 4055.7273
$ws.Range("K132").Value = 3002.4
$ws.Range("L132").Value = 12167.1819
$ws.Range("M132").Value = -472.3999999999996
$ws.Range("N132").Value = -17227.1819
